$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "teenAge" and fill it with the new age-grouping table ---
$teenAge = $wb.Worksheets.Item("Sheet3")
$teenAge.Name = "teenAge"

$teenAge.Range("A1").Value = "lAge"
$teenAge.Range("B1").Value = "uAge"
$teenAge.Range("A2").Value = 0
$teenAge.Range("B2").Value = 17
$teenAge.Range("A3").Value = 18
$teenAge.Range("B3").Value = 64
$teenAge.Range("A4").Value = 65
$teenAge.Range("B4").Value = 999

$teenAge.Range("A10").Select()

# --- Duplicate teenAge into a new "drinkAge" sheet right after it ---
$teenAge.Copy($null, $teenAge)
$drinkAge = $wb.Worksheets.Item($wb.Worksheets.Count)
$drinkAge.Name = "drinkAge"
$drinkAge.Activate()
$drinkAge.Range("E8").Select()

# teenAge (but not drinkAge) gets an explicit portrait page setup
$teenAge.PageSetup.Orientation = 1

# --- Update selection on the "data" sheet ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("A1:B4").Select()

# --- Re-activate drinkAge so it ends up the active/visible tab ---
$drinkAge.Activate()
